$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 (PORCELANATO)
$ws.Range("D3").Value = 552.4299999999999
$ws.Range("E3").Value = 13170.91
$ws.Range("F3").Value = 0.04025477762702082

# Row 4 (TOTAL)
$ws.Range("D4").Value = 1478.34
$ws.Range("E4").Value = 12245
$ws.Range("F4").Value = 0.1077245043845011
